$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite the DNA-part transfer rows 10-21 (UID 9-20): the previous two
# rows are updated and 10 new rows are appended, covering the full source
# well x destination well combinations.
$data = @(
    @(9,  "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A1", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A3", 250, "BBa_promoter"),
    @(10, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A2", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A3", 250, "BBa_rbs"),
    @(11, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A3", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A3", 250, "BBa_CDSrestrict"),
    @(12, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A4", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A3", 250, "BBa_term"),
    @(13, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A1", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A4", 250, "BBa_promoter"),
    @(14, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A2", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A4", 250, "BBa_rbs"),
    @(15, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A3", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A4", 250, "BBa_CDSrestrict"),
    @(16, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A4", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A4", 250, "BBa_term"),
    @(17, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A5", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A1", 250, "pTU1-A-RFP"),
    @(18, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A5", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A2", 500, "pTU1-A-RFP"),
    @(19, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A6", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A3", 250, "pTU1-B-RFP"),
    @(20, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A6", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A4", 500, "pTU1-B-RFP")
)

$rowCount = $data.Count
$colCount = $data[0].Count

$arr = New-Object 'object[,]' $rowCount,$colCount
for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $arr[$r, $c] = $data[$r][$c]
    }
}

$startRow = 10
$endRow = $startRow + $rowCount - 1
$ws.Range("A$startRow`:I$endRow").Value = $arr
